$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "date" cell (B5) so the new cells reuse the existing shared
# string "2019-01-04" instead of Excel auto-converting typed text to a
# date value/format.
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4104)
$ws.Range("B5").Copy()
$ws.Range("B7").PasteSpecial(-4104)
$ws.Range("B5").Copy()
$ws.Range("B8").PasteSpecial(-4104)

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 3).Value = 117.88

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 3).Value = 577.85

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 3).Value = 200.4
